# Remove the "GAPTIN 400 MG 30 CAPS." line (row 41) from the day-sale report
# and refresh the generated-at timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole row for GAPTIN 400 MG 30 CAPS. - everything below shifts up.
$ws.Rows("41:41").Delete()

# Column A just holds the sequential item number (1..93 after the delete); the
# native row-delete shifted those literal numbers up along with everything
# else, so renumber them back to a plain 1..N sequence.
for ($r = 41; $r -le 99; $r++) {
    $ws.Range("A$r").Value = $r - 6
}

# The grand-total cell is a hard-coded number (not a formula), so update it
# by hand now that GAPTIN's price (51.48) has dropped out of the report.
$ws.Range("P100").Value = 4946.6

# Refresh the "generated at" timestamp footer to the new save time.
$ws.Range("A101").Value = "Tuesday, 17 June, 2025 8:45 PM"
